# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" worksheets to match the refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 2872
    4  = 98
    5  = 6695
    6  = 1605
    7  = 17
    9  = 52
    11 = 21
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
